# Auto-generated Excel COM-interop script to apply F-column updates
# across worksheets "展览" (sheet1), "演出" (sheet2), "全部类型" (sheet4).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 632
$ws.Range("F3").Value = 663
$ws.Range("F4").Value = 926
$ws.Range("F5").Value = 687
$ws.Range("F6").Value = 825
$ws.Range("F7").Value = 384
$ws.Range("F8").Value = 587
$ws.Range("F9").Value = 122
$ws.Range("F10").Value = 1186
$ws.Range("F11").Value = 615
$ws.Range("F12").Value = 363
$ws.Range("F13").Value = 490
$ws.Range("F15").Value = 283
$ws.Range("F16").Value = 327
$ws.Range("F17").Value = 49
$ws.Range("F19").Value = 541
$ws.Range("F20").Value = 55
$ws.Range("F21").Value = 558
$ws.Range("F22").Value = 23
$ws.Range("F23").Value = 636

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 309
$ws.Range("F6").Value = 18
$ws.Range("F9").Value = 214
$ws.Range("F10").Value = 46
$ws.Range("F11").Value = 20
$ws.Range("F13").Value = 55

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 632
$ws.Range("F6").Value = 309
$ws.Range("F7").Value = 663
$ws.Range("F8").Value = 926
$ws.Range("F9").Value = 687
$ws.Range("F10").Value = 825
$ws.Range("F11").Value = 384
$ws.Range("F12").Value = 587
$ws.Range("F13").Value = 122
$ws.Range("F14").Value = 1186
$ws.Range("F15").Value = 615
$ws.Range("F17").Value = 18
$ws.Range("F18").Value = 363
$ws.Range("F19").Value = 490
$ws.Range("F22").Value = 283
$ws.Range("F24").Value = 327
$ws.Range("F25").Value = 49
$ws.Range("F27").Value = 214
$ws.Range("F28").Value = 46
$ws.Range("F29").Value = 541
$ws.Range("F30").Value = 20
$ws.Range("F32").Value = 55
$ws.Range("F33").Value = 55
$ws.Range("F34").Value = 558
$ws.Range("F35").Value = 23
$ws.Range("F36").Value = 636

